$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

$ws.Range("P2").Value = "['[1,1,2]', '[0,2,1]', '[2,1,0]', '[1,2,0]', '[0,0,0]', '[1,1,1]', '[0,1,1]', '[2,0,1]', '[1,0,2]', '[2,2,2]']"
$ws.Range("P3").Value = "['[1,1,2]', '[0,2,1]', '[2,1,0]', '[1,2,0]', '[0,0,0]', '[1,1,1]', '[0,1,1]', '[2,0,1]', '[1,0,2]', '[2,2,2]']"

$ws.Range("P4").Value = "['[1,2,0]', '[1,0,2]', '[0,0,0]', '[1,1,1]', '[2,0,1]', '[2,2,2]', '[2,1,0]', '[0,1,1]', '[1,1,2]', '[0,2,1]']"
$ws.Range("P5").Value = "['[1,2,0]', '[1,0,2]', '[0,0,0]', '[1,1,1]', '[2,0,1]', '[2,2,2]', '[2,1,0]', '[0,1,1]', '[1,1,2]', '[0,2,1]']"

$ws.Range("P6").Value = "['[2,1,0]', '[1,2,0]', '[1,1,1]', '[0,1,1]', '[1,0,2]', '[0,2,1]', '[2,2,2]', '[1,1,2]', '[2,0,1]', '[0,0,0]']"
$ws.Range("P7").Value = "['[2,1,0]', '[1,2,0]', '[1,1,1]', '[0,1,1]', '[1,0,2]', '[0,2,1]', '[2,2,2]', '[1,1,2]', '[2,0,1]', '[0,0,0]']"

$ws.Range("P8").Value = "['[2,0,1]', '[1,1,1]', '[0,0,0]', '[2,2,2]', '[2,1,0]', '[1,1,2]', '[1,2,0]', '[0,1,1]', '[1,0,2]', '[0,2,1]']"
$ws.Range("P9").Value = "['[2,0,1]', '[1,1,1]', '[0,0,0]', '[2,2,2]', '[2,1,0]', '[1,1,2]', '[1,2,0]', '[0,1,1]', '[1,0,2]', '[0,2,1]']"

$ws.Range("P10").Value = "['[1,2,0]', '[2,1,0]', '[0,0,0]', '[1,1,1]', '[0,1,1]', '[2,2,2]', '[1,0,2]', '[1,1,2]', '[0,2,1]', '[2,0,1]']"
$ws.Range("P11").Value = "['[1,2,0]', '[2,1,0]', '[0,0,0]', '[1,1,1]', '[0,1,1]', '[2,2,2]', '[1,0,2]', '[1,1,2]', '[0,2,1]', '[2,0,1]']"

$ws.Range("P12").Value = "['[0,0,0]', '[0,2,1]', '[1,1,2]', '[1,0,2]', '[1,2,0]', '[2,2,2]', '[2,0,1]', '[2,1,0]', '[1,1,1]', '[0,1,1]']"
$ws.Range("P13").Value = "['[0,0,0]', '[0,2,1]', '[1,1,2]', '[1,0,2]', '[1,2,0]', '[2,2,2]', '[2,0,1]', '[2,1,0]', '[1,1,1]', '[0,1,1]']"

$ws.Range("P14").Value = "['[1,1,2]', '[2,0,1]', '[2,2,2]', '[1,0,2]', '[1,1,1]', '[0,1,1]', '[2,1,0]', '[0,2,1]', '[0,0,0]', '[1,2,0]']"
$ws.Range("P15").Value = "['[1,1,2]', '[2,0,1]', '[2,2,2]', '[1,0,2]', '[1,1,1]', '[0,1,1]', '[2,1,0]', '[0,2,1]', '[0,0,0]', '[1,2,0]']"

$ws.Range("P16").Value = "['[0,1,1]', '[2,2,2]', '[2,1,0]', '[0,0,0]', '[1,0,2]', '[0,2,1]', '[1,1,2]', '[2,0,1]', '[1,1,1]', '[1,2,0]']"
$ws.Range("P17").Value = "['[0,1,1]', '[2,2,2]', '[2,1,0]', '[0,0,0]', '[1,0,2]', '[0,2,1]', '[1,1,2]', '[2,0,1]', '[1,1,1]', '[1,2,0]']"

$ws.Range("P18").Value = "['[0,0,0]', '[1,0,2]', '[1,2,0]', '[0,1,1]', '[1,1,1]', '[2,1,0]', '[2,0,1]', '[1,1,2]', '[2,2,2]', '[0,2,1]']"
$ws.Range("P19").Value = "['[0,0,0]', '[1,0,2]', '[1,2,0]', '[0,1,1]', '[1,1,1]', '[2,1,0]', '[2,0,1]', '[1,1,2]', '[2,2,2]', '[0,2,1]']"

$ws.Range("P20").Value = "['[2,2,2]', '[1,1,1]', '[2,1,0]', '[0,1,1]', '[1,2,0]', '[0,2,1]', '[1,1,2]', '[1,0,2]', '[0,0,0]', '[2,0,1]']"
$ws.Range("P21").Value = "['[2,2,2]', '[1,1,1]', '[2,1,0]', '[0,1,1]', '[1,2,0]', '[0,2,1]', '[1,1,2]', '[1,0,2]', '[0,0,0]', '[2,0,1]']"
